$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1) Update the "总计" (summary) sheet: insert a new row for 2022-Q3 above the
#    existing 2022-Q2 row, shifting the old row down.
# ---------------------------------------------------------------------------
$summary = $wb.Worksheets.Item("总计")

$summary.Rows.Item(2).Insert()

# Restore the index-column style (s=2) on the newly inserted A2, which Insert()
# leaves blank/unstyled, by copying the format from A3 (the shifted-down row
# that kept the original style).
$summary.Range("A3").Copy()
$summary.Range("A2").PasteSpecial(-4122)

# New 2022-Q3 summary row
$summary.Cells(2, 1).Value = 0
$summary.Cells(2, 2).Value = "2022-Q3"
$summary.Cells(2, 3).Value = 4
$summary.Cells(2, 4).Value = 0.08

# Renumber / rewrite the shifted-down 2022-Q2 row explicitly
$summary.Cells(3, 1).Value = 1
$summary.Cells(3, 2).Value = "2022-Q2"
$summary.Cells(3, 3).Value = 1
$summary.Cells(3, 4).Value = 0

# The data cells (B:D) of the newly inserted row picked up stray formatting
# from Insert()'s "copy format from row above" behaviour; the source row has
# no explicit style on these cells, so strip it back off.
$summary.Range("B2:D2").ClearFormats()

# ---------------------------------------------------------------------------
# 2) Insert a brand-new "2022-Q3" fund-holdings sheet, positioned between
#    "总计" and "2022-Q2".
# ---------------------------------------------------------------------------
$q2ref = $wb.Worksheets.Item("2022-Q2")
$q3 = $wb.Worksheets.Add($q2ref)
$q3.Name = "2022-Q3"

# Match styling (s=2) used elsewhere in this workbook for header rows / the
# numeric index column, by copying format from the "总计" sheet.
$summary.Range("B1:D1").Copy()
$q3.Range("B1:H1").PasteSpecial(-4122)
$summary.Range("A2").Copy()
$q3.Range("A2:A5").PasteSpecial(-4122)

$headers = @("基金代码", "基金名称", "基金规模", "股票总仓位", "仓位占比", "持有市值(亿元)", "仓位排名")
for ($i = 0; $i -lt $headers.Length; $i++) {
    $q3.Cells(1, $i + 2).Value = $headers[$i]
}

# Force columns B:G to store as text (matching the source data's inlineStr
# cells, e.g. fund codes with leading zeros like "015143") instead of being
# auto-coerced to numbers.
$q3.Range("B2:G5").NumberFormat = "@"

$funds = @(
    @(0, "015143", "中欧智能制造混合A", "0.96", "84.58", "3.31", "0.0318", 9),
    @(1, "690003", "民生加银精选混合", "0.51", "89.96", "5.63", "0.0287", 7),
    @(2, "015144", "中欧智能制造混合C", "0.51", "84.58", "3.31", "0.0169", 9),
    @(3, "005706", "兴业龙腾双益平衡混合", "1.81", "32.02", "0.16", "0.0029", 10)
)

for ($r = 0; $r -lt $funds.Length; $r++) {
    $row = $funds[$r]
    $excelRow = $r + 2
    $q3.Cells($excelRow, 1).Value = $row[0]
    $q3.Cells($excelRow, 2).Value = $row[1]
    $q3.Cells($excelRow, 3).Value = $row[2]
    $q3.Cells($excelRow, 4).Value = $row[3]
    $q3.Cells($excelRow, 5).Value = $row[4]
    $q3.Cells($excelRow, 6).Value = $row[5]
    $q3.Cells($excelRow, 7).Value = $row[6]
    $q3.Cells($excelRow, 8).Value = $row[7]
}

# Drop the temporary "@" text number-format now that the values are stored as
# text, so the data cells end up unstyled (matching the source).
$q3.Range("B2:G5").ClearFormats()

# ---------------------------------------------------------------------------
# 3) Restore "2022-Q2" as the selected/active sheet (it was the active sheet
#    before this edit, and adding a new sheet shifts selection to it).
# ---------------------------------------------------------------------------
$wb.Worksheets.Item("2022-Q2").Select()
